# Auto-applied scheduled market-data refresh for Siren_Profits workbook.
# Updates currentAveragePrice(NQ/HQ), LevePrice(NQ/HQ) and LeveProfit(NQ/HQ)
# columns (H-N) on affected leve rows across all job sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 55
$ws.Cells.Item(55, 8).Value = 96.59999999999999
$ws.Cells.Item(55, 9).Value = 114.75
$ws.Cells.Item(55, 10).Value = 24
$ws.Cells.Item(55, 11).Value = 114.75
$ws.Cells.Item(55, 12).Value = 24
$ws.Cells.Item(55, 13).Value = 99.25
$ws.Cells.Item(55, 14).Value = -452
# Row 74
$ws.Cells.Item(74, 8).Value = 5333.05
$ws.Cells.Item(74, 9).Value = 4946.3335
$ws.Cells.Item(74, 10).Value = 5498.7856
$ws.Cells.Item(74, 11).Value = 4946.3335
$ws.Cells.Item(74, 12).Value = 5498.7856
$ws.Cells.Item(74, 13).Value = -4010.3335
$ws.Cells.Item(74, 14).Value = -7370.7856
# Row 77
$ws.Cells.Item(77, 8).Value = 5333.05
$ws.Cells.Item(77, 9).Value = 4946.3335
$ws.Cells.Item(77, 10).Value = 5498.7856
$ws.Cells.Item(77, 11).Value = 24731.6675
$ws.Cells.Item(77, 12).Value = 27493.928
$ws.Cells.Item(77, 13).Value = -20051.6675
$ws.Cells.Item(77, 14).Value = -36853.928
# Row 111
$ws.Cells.Item(111, 8).Value = 3467
$ws.Cells.Item(111, 9).Value = 3111.6
$ws.Cells.Item(111, 10).Value = 3763.1667
$ws.Cells.Item(111, 11).Value = 9334.799999999999
$ws.Cells.Item(111, 12).Value = 11289.5001
$ws.Cells.Item(111, 13).Value = -6267.799999999999
$ws.Cells.Item(111, 14).Value = -17423.5001
# Row 114
$ws.Cells.Item(114, 8).Value = 722000
$ws.Cells.Item(114, 10).Value = 722000
$ws.Cells.Item(114, 12).Value = 722000
$ws.Cells.Item(114, 14).Value = -730678
# Row 125
$ws.Cells.Item(125, 8).Value = 2248.6667
$ws.Cells.Item(125, 9).Value = 1246
$ws.Cells.Item(125, 11).Value = 11214
$ws.Cells.Item(125, 13).Value = -8754
# Row 138
$ws.Cells.Item(138, 8).Value = 3646.913
$ws.Cells.Item(138, 9).Value = 1590.25
$ws.Cells.Item(138, 11).Value = 4770.75
$ws.Cells.Item(138, 13).Value = 369.25

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 37
$ws.Cells.Item(37, 8).Value = 30017
$ws.Cells.Item(37, 10).Value = 30000
$ws.Cells.Item(37, 12).Value = 30000
$ws.Cells.Item(37, 14).Value = -30546
# Row 74
$ws.Cells.Item(74, 8).Value = 4440.2856
$ws.Cells.Item(74, 9).Value = 2569.9473
$ws.Cells.Item(74, 10).Value = 8388.777
$ws.Cells.Item(74, 11).Value = 2569.9473
$ws.Cells.Item(74, 12).Value = 8388.777
$ws.Cells.Item(74, 13).Value = -1695.9473
$ws.Cells.Item(74, 14).Value = -10136.777
# Row 77
$ws.Cells.Item(77, 8).Value = 4440.2856
$ws.Cells.Item(77, 9).Value = 2569.9473
$ws.Cells.Item(77, 10).Value = 8388.777
$ws.Cells.Item(77, 11).Value = 12849.7365
$ws.Cells.Item(77, 12).Value = 41943.885
$ws.Cells.Item(77, 13).Value = -8481.736499999999
$ws.Cells.Item(77, 14).Value = -50679.885
# Row 122
$ws.Cells.Item(122, 8).Value = 363429.34
$ws.Cells.Item(122, 9).Value = 2675.3044
$ws.Cells.Item(122, 11).Value = 8025.9132
$ws.Cells.Item(122, 13).Value = -5575.9132
# Row 132
$ws.Cells.Item(132, 8).Value = 2664.68
$ws.Cells.Item(132, 9).Value = 1577.2632
$ws.Cells.Item(132, 11).Value = 4731.7896
$ws.Cells.Item(132, 13).Value = -2201.7896
# Row 133
$ws.Cells.Item(133, 8).Value = 63628.2
$ws.Cells.Item(133, 10).Value = 63628.2
$ws.Cells.Item(133, 12).Value = 63628.2
$ws.Cells.Item(133, 14).Value = -68688.2

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Cells.Item(105, 8).Value = 2787.2693
$ws.Cells.Item(105, 9).Value = 1182.6316
$ws.Cells.Item(105, 10).Value = 7142.7144
$ws.Cells.Item(105, 11).Value = 1182.6316
$ws.Cells.Item(105, 12).Value = 7142.7144
$ws.Cells.Item(105, 13).Value = 564.3684000000001
$ws.Cells.Item(105, 14).Value = -10636.7144

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Cells.Item(16, 8).Value = 1226.6666
$ws.Cells.Item(16, 9).Value = 1226.6666
$ws.Cells.Item(16, 11).Value = 1226.6666
$ws.Cells.Item(16, 13).Value = -939.6666
# Row 22
$ws.Cells.Item(22, 8).Value = 532.3333
$ws.Cells.Item(22, 9).Value = 438.8
$ws.Cells.Item(22, 10).Value = 1000
$ws.Cells.Item(22, 11).Value = 438.8
$ws.Cells.Item(22, 12).Value = 1000
$ws.Cells.Item(22, 13).Value = -88.80000000000001
$ws.Cells.Item(22, 14).Value = -1700
# Row 31
$ws.Cells.Item(31, 8).Value = 2848.8
$ws.Cells.Item(31, 9).Value = 1731.8667
$ws.Cells.Item(31, 10).Value = 6199.6
$ws.Cells.Item(31, 11).Value = 1731.8667
$ws.Cells.Item(31, 12).Value = 6199.6
$ws.Cells.Item(31, 13).Value = -1436.8667
$ws.Cells.Item(31, 14).Value = -6789.6
# Row 34
$ws.Cells.Item(34, 8).Value = 2848.8
$ws.Cells.Item(34, 9).Value = 1731.8667
$ws.Cells.Item(34, 10).Value = 6199.6
$ws.Cells.Item(34, 11).Value = 1731.8667
$ws.Cells.Item(34, 12).Value = 6199.6
$ws.Cells.Item(34, 13).Value = -1529.8667
$ws.Cells.Item(34, 14).Value = -6603.6
# Row 99
$ws.Cells.Item(99, 8).Value = 369805.84
$ws.Cells.Item(99, 10).Value = 6999
$ws.Cells.Item(99, 12).Value = 6999
$ws.Cells.Item(99, 14).Value = -9995
# Row 113
$ws.Cells.Item(113, 8).Value = 1226.6666
$ws.Cells.Item(113, 9).Value = 1226.6666
$ws.Cells.Item(113, 11).Value = 1226.6666
$ws.Cells.Item(113, 13).Value = 943.3334
# Row 126
$ws.Cells.Item(126, 8).Value = 369805.84
$ws.Cells.Item(126, 10).Value = 6999
$ws.Cells.Item(126, 12).Value = 20997
$ws.Cells.Item(126, 14).Value = -25937
# Row 138
$ws.Cells.Item(138, 8).Value = 84089.37
$ws.Cells.Item(138, 10).Value = 84089.37
$ws.Cells.Item(138, 12).Value = 84089.37
$ws.Cells.Item(138, 14).Value = -94369.37

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 20
$ws.Cells.Item(20, 8).Value = 150
$ws.Cells.Item(20, 9).Value = 0
$ws.Cells.Item(20, 10).Value = 150
$ws.Cells.Item(20, 11).Value = 0
$ws.Cells.Item(20, 12).Value = 450
$ws.Cells.Item(20, 13).Value = $null
$ws.Cells.Item(20, 14).Value = -904
# Row 38
$ws.Cells.Item(38, 8).Value = 1209.5938
$ws.Cells.Item(38, 9).Value = 322.5
$ws.Cells.Item(38, 11).Value = 967.5
$ws.Cells.Item(38, 13).Value = -620.5

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 48
$ws.Cells.Item(48, 8).Value = 25000
$ws.Cells.Item(48, 10).Value = 25000
$ws.Cells.Item(48, 12).Value = 25000
$ws.Cells.Item(48, 14).Value = -25970
# Row 123
$ws.Cells.Item(123, 8).Value = 29750
$ws.Cells.Item(123, 10).Value = 29750
$ws.Cells.Item(123, 12).Value = 29750
$ws.Cells.Item(123, 14).Value = -34650
# Row 126
$ws.Cells.Item(126, 8).Value = 15828.346
$ws.Cells.Item(126, 10).Value = 12829.059
$ws.Cells.Item(126, 12).Value = 38487.177
$ws.Cells.Item(126, 14).Value = -43427.177

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Cells.Item(22, 8).Value = 8357.809999999999
$ws.Cells.Item(22, 9).Value = 10638.143
$ws.Cells.Item(22, 10).Value = 3797.1428
$ws.Cells.Item(22, 11).Value = 10638.143
$ws.Cells.Item(22, 12).Value = 3797.1428
$ws.Cells.Item(22, 13).Value = -10343.143
$ws.Cells.Item(22, 14).Value = -4387.1428
# Row 27
$ws.Cells.Item(27, 8).Value = 8357.809999999999
$ws.Cells.Item(27, 9).Value = 10638.143
$ws.Cells.Item(27, 10).Value = 3797.1428
$ws.Cells.Item(27, 11).Value = 10638.143
$ws.Cells.Item(27, 12).Value = 3797.1428
$ws.Cells.Item(27, 13).Value = -10531.143
$ws.Cells.Item(27, 14).Value = -4011.1428
# Row 46
$ws.Cells.Item(46, 8).Value = 1688
$ws.Cells.Item(46, 9).Value = 878.4286
$ws.Cells.Item(46, 10).Value = 2632.5
$ws.Cells.Item(46, 11).Value = 878.4286
$ws.Cells.Item(46, 12).Value = 2632.5
$ws.Cells.Item(46, 13).Value = -690.4286
$ws.Cells.Item(46, 14).Value = -3008.5
# Row 118
$ws.Cells.Item(118, 8).Value = 0
$ws.Cells.Item(118, 10).Value = 0
$ws.Cells.Item(118, 12).Value = 0
$ws.Cells.Item(118, 14).Value = $null
# Row 122
$ws.Cells.Item(122, 8).Value = 3753.875
$ws.Cells.Item(122, 9).Value = 3262.087
$ws.Cells.Item(122, 11).Value = 9786.261
$ws.Cells.Item(122, 13).Value = -7336.261
# Row 132
$ws.Cells.Item(132, 8).Value = 416667.62
$ws.Cells.Item(132, 9).Value = 678897.0600000001
$ws.Cells.Item(132, 11).Value = 2036691.18
$ws.Cells.Item(132, 13).Value = -2034161.18

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Cells.Item(107, 8).Value = 29342.318
$ws.Cells.Item(107, 9).Value = 2363
$ws.Cells.Item(107, 11).Value = 7089
$ws.Cells.Item(107, 13).Value = -5169
# Row 122
$ws.Cells.Item(122, 8).Value = 5126.1924
$ws.Cells.Item(122, 9).Value = 3915.25
$ws.Cells.Item(122, 10).Value = 6164.143
$ws.Cells.Item(122, 11).Value = 11745.75
$ws.Cells.Item(122, 12).Value = 18492.429
$ws.Cells.Item(122, 13).Value = -9295.75
$ws.Cells.Item(122, 14).Value = -23392.429
# Row 123
$ws.Cells.Item(123, 8).Value = 30389.4
$ws.Cells.Item(123, 10).Value = 30385.5
$ws.Cells.Item(123, 12).Value = 30385.5
$ws.Cells.Item(123, 14).Value = -40185.5
